$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 data (append before formatting so the used range grows first)
$ws.Range("A8").Value = "SMIJO"
$ws.Range("B8").Value = "cleverladyiam@gmail.com"
$ws.Range("C8").Value = "ALCHP-INS-202514445"
$ws.Range("D8").Value = "Joby Chirayath House, Pazhuvil PO, Pin 680564, Thrissur"
$ws.Range("E8").Value = 7356140066
$ws.Range("F8").Value = 20025
$ws.Range("G8").Value = "CROWN-114445"

# Column D (Address) + Column E (Number) corrections for existing rows
$ws.Range("D2").Value = "Karunya, Ambancode, Peyad PO, Trivandrum, Kerala, INDIA, Pincode: 695573"
$ws.Range("E2").Value = 9656005657

$ws.Range("D3").Value = "Revathy House, Ayyankavu Temple Road, Irinjalakuda, Thrissur, Pincode: 680121"
$ws.Range("E3").Value = 8075937035

$ws.Range("D4").Value = "Revathy House, Ayyankavu Temple Road, Irinjalakuda, Thrissur, Pincode: 680121"
$ws.Range("E4").Value = 9496143678

$ws.Range("D5").Value = "Revathy House, Ayyankavu Temple Road, Irinjalakuda, Thrissur, Pincode: 680121"
$ws.Range("E5").Value = 9496143678

$ws.Range("D6").Value = "Thekkekara House,Chettiparambu, Irinjalakuda, Thrissur District, Kerala, India "
$ws.Range("E6").Value = 9946580236

$ws.Range("D7").Value = "Thekkekara House,Chettiparambu, Irinjalakuda, Thrissur District, Kerala, India "
$ws.Range("E7").Value = 9946580235

# Apply the "General" number format across the full used range, which mints
# the new numFmt/cellXf (id 164) and stamps s="1" on every cell
$ws.Range("A1:G8").NumberFormat = "General"
